# Updated cryptos list on Thu Mar 28 07:59:45 UTC 2024 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures scraped from coinranking.com;
# three coin pairs also swapped rank order (rows 35/36, 46/47, 50/51).
#
# Price cells are plain text in the sheet (e.g. thousand-separator dots like
# '70.352.74' wouldn't round-trip as numbers). A leading apostrophe forces
# Excel to keep any numeric-looking price string (e.g. '590.04') as text too,
# matching the existing column formatting without touching cell styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '70.352.74'
$ws.Cells.Item(2, 5).Value = '  +0.90%  '

$ws.Cells.Item(3, 4).Value = '3.570.75'
$ws.Cells.Item(3, 5).Value = '  +0.19%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).Value = '''590.04'
$ws.Cells.Item(5, 5).Value = '  +2.65%  '

$ws.Cells.Item(6, 4).Value = '''187.35'
$ws.Cells.Item(6, 5).Value = '  +0.54%  '

$ws.Cells.Item(7, 4).Value = '3.560.49'
$ws.Cells.Item(7, 5).Value = '  -0.08%  '

$ws.Cells.Item(8, 4).Value = '''0.622'
$ws.Cells.Item(8, 5).Value = '  +0.25%  '

$ws.Cells.Item(9, 5).Value = '  +0.01%  '

$ws.Cells.Item(10, 5).Value = '  +8.97%  '

$ws.Cells.Item(11, 4).Value = '''0.650'
$ws.Cells.Item(11, 5).Value = '  +0.44%  '

$ws.Cells.Item(12, 4).Value = '''54.84'
$ws.Cells.Item(12, 5).Value = '  -0.11%  '

$ws.Cells.Item(13, 5).Value = '  +1.86%  '

$ws.Cells.Item(14, 5).Value = '  +0.84%  '

$ws.Cells.Item(15, 4).Value = '4.136.30'
$ws.Cells.Item(15, 5).Value = '  -0.09%  '

$ws.Cells.Item(16, 4).Value = '''19.49'
$ws.Cells.Item(16, 5).Value = '  -0.19%  '

$ws.Cells.Item(17, 4).Value = '70.357.11'
$ws.Cells.Item(17, 5).Value = '  +0.93%  '

$ws.Cells.Item(18, 4).Value = '3.556.05'
$ws.Cells.Item(18, 5).Value = '  -0.39%  '

$ws.Cells.Item(19, 4).Value = '''12.49'
$ws.Cells.Item(19, 5).Value = '  +0.23%  '

$ws.Cells.Item(20, 5).Value = '  -0.66%  '

$ws.Cells.Item(21, 4).Value = '''557.65'
$ws.Cells.Item(21, 5).Value = '  +14.41%  '

$ws.Cells.Item(22, 5).Value = '  -0.29%  '

$ws.Cells.Item(23, 4).Value = '''17.97'
$ws.Cells.Item(23, 5).Value = '  -8.25%  '

$ws.Cells.Item(24, 4).Value = '''4.69'
$ws.Cells.Item(24, 5).Value = '  +9.11%  '

$ws.Cells.Item(25, 5).Value = '  +0.26%  '

$ws.Cells.Item(26, 4).Value = '''96.29'
$ws.Cells.Item(26, 5).Value = '  +0.75%  '

$ws.Cells.Item(27, 4).Value = '''11.53'
$ws.Cells.Item(27, 5).Value = '  +4.84%  '

$ws.Cells.Item(28, 5).Value = '  +1.81%  '

$ws.Cells.Item(29, 4).Value = '''9.18'
$ws.Cells.Item(29, 5).Value = '  -0.86%  '

$ws.Cells.Item(30, 4).Value = '''32.34'
$ws.Cells.Item(30, 5).Value = '  +2.23%  '

$ws.Cells.Item(31, 4).Value = '''7.36'
$ws.Cells.Item(31, 5).Value = '  -1.66%  '

$ws.Cells.Item(32, 4).Value = '''12.57'
$ws.Cells.Item(32, 5).Value = '  +4.93%  '

$ws.Cells.Item(33, 4).Value = '''65.22'
$ws.Cells.Item(33, 5).Value = '  -1.23%  '

$ws.Cells.Item(34, 5).Value = '  -0.39%  '

$ws.Cells.Item(35, 2).Value = 'Bittensor'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(35, 4).Value = '''557.46'
$ws.Cells.Item(35, 5).Value = '  -1.98%  '

$ws.Cells.Item(36, 2).Value = 'Fetch.AI'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(36, 4).Value = '''3.26'
$ws.Cells.Item(36, 5).Value = '  +6.47%  '

$ws.Cells.Item(37, 4).Value = '''0.418'
$ws.Cells.Item(37, 5).Value = '  +7.35%  '

$ws.Cells.Item(38, 4).Value = '''38.27'
$ws.Cells.Item(38, 5).Value = '  +0.24%  '

$ws.Cells.Item(39, 4).Value = '''0.999'

$ws.Cells.Item(40, 4).Value = '0.0₃0772'
$ws.Cells.Item(40, 5).Value = '  -2.76%  '

$ws.Cells.Item(41, 5).Value = '  -0.28%  '

$ws.Cells.Item(42, 4).Value = '3.373.36'
$ws.Cells.Item(42, 5).Value = '  +3.51%  '

$ws.Cells.Item(43, 4).Value = '''3.12'
$ws.Cells.Item(43, 5).Value = '  -3.91%  '

$ws.Cells.Item(44, 5).Value = '  -2.82%  '

$ws.Cells.Item(45, 5).Value = '  +3.01%  '

$ws.Cells.Item(46, 2).Value = 'ThetaToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(46, 4).Value = '''2.99'
$ws.Cells.Item(46, 5).Value = '  +0.40%  '

$ws.Cells.Item(47, 2).Value = 'VeChain'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(47, 4).Value = '''0.0447'
$ws.Cells.Item(47, 5).Value = '  +3.02%  '

$ws.Cells.Item(48, 4).Value = '''9.19'
$ws.Cells.Item(48, 5).Value = '  -4.32%  '

$ws.Cells.Item(49, 5).Value = '  +0.87%  '

$ws.Cells.Item(50, 2).Value = 'OceanProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Cells.Item(50, 4).Value = '''1.50'
$ws.Cells.Item(50, 5).Value = '  +24.95%  '

$ws.Cells.Item(51, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(51, 4).Value = '''1.00'
$ws.Cells.Item(51, 5).Value = '  +0.03%  '
